$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 3-5 entirely (they are removed from the sheet)
$ws.Range("A3:B5").EntireRow.Delete() | Out-Null

# Update row 1 (header) values
$ws.Range("A1").Value = "aca"
$ws.Range("B1").Value = "aja"
$ws.Range("C1").Value = "aoa"

# Copy the header style (s="1") to the new C1 cell so it matches A1/B1
$ws.Range("A1").Copy() | Out-Null
$ws.Range("C1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# Update row 2 values
$ws.Range("A2").Value = "aca"
$ws.Range("B2").Value = "aja"
$ws.Range("C2").Value = "aoa"
